# Update "想去人数" (interest counter, column F) figures on each sheet
# to match the regenerated site data (gh-pages output at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 14
$ws.Cells.Item(3, 6).Value = 2723
$ws.Cells.Item(4, 6).Value = 1064
$ws.Cells.Item(5, 6).Value = 19591
$ws.Cells.Item(6, 6).Value = 75
$ws.Cells.Item(7, 6).Value = 2211
$ws.Cells.Item(9, 6).Value = 612
$ws.Cells.Item(10, 6).Value = 432
$ws.Cells.Item(12, 6).Value = 236
$ws.Cells.Item(15, 6).Value = 363
$ws.Cells.Item(18, 6).Value = 169
$ws.Cells.Item(19, 6).Value = 187
$ws.Cells.Item(20, 6).Value = 19
$ws.Cells.Item(22, 6).Value = 99
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 193
$ws.Cells.Item(3, 6).Value = 36
$ws.Cells.Item(5, 6).Value = 14
$ws.Cells.Item(7, 6).Value = 286
$ws.Cells.Item(8, 6).Value = 129
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 5998
$ws.Cells.Item(3, 6).Value = 640
$ws.Cells.Item(4, 6).Value = 590
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 5998
$ws.Cells.Item(3, 6).Value = 640
$ws.Cells.Item(4, 6).Value = 590
$ws.Cells.Item(5, 6).Value = 193
$ws.Cells.Item(6, 6).Value = 14
$ws.Cells.Item(7, 6).Value = 36
$ws.Cells.Item(8, 6).Value = 2723
$ws.Cells.Item(9, 6).Value = 1064
$ws.Cells.Item(10, 6).Value = 19591
$ws.Cells.Item(12, 6).Value = 14
$ws.Cells.Item(13, 6).Value = 75
$ws.Cells.Item(15, 6).Value = 286
$ws.Cells.Item(16, 6).Value = 2211
$ws.Cells.Item(18, 6).Value = 129
$ws.Cells.Item(19, 6).Value = 612
$ws.Cells.Item(20, 6).Value = 432
$ws.Cells.Item(22, 6).Value = 236
$ws.Cells.Item(28, 6).Value = 363
$ws.Cells.Item(33, 6).Value = 169
$ws.Cells.Item(35, 6).Value = 187
$ws.Cells.Item(38, 6).Value = 19
$ws.Cells.Item(47, 6).Value = 99
